# Implemented static division by zero checks
# Adds two new "Todo" rows documenting optimisation/analysis checks that
# are still outstanding: "Dead code" and "Loop Unrolling".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27 - Dead code elimination check
$ws.Range("A27").Value = "?"
$ws.Range("B27").Value = "Me"
$ws.Range("C27").Value = "Dead code"
$ws.Range("D27").Value = "Y"

# Row 28 - Loop unrolling check
$ws.Range("A28").Value = "?"
$ws.Range("B28").Value = "Me"
$ws.Range("C28").Value = "Loop Unrolling"
$ws.Range("D28").Value = "Y"

# Move the active selection to the newly added last row, matching Excel's
# behaviour after entering data in sequence down column C.
$ws.Range("C28").Select()
